$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the single remaining worker row (row 16) with the data that used
# to live in row 18 (LEIDIS PAOLA BLANCO DE LA ROSA).
$ws.Range("C16").Value = "1047452509"
$ws.Range("D16").Value = "LEIDIS PAOLA BLANCO DE LA ROSA"
$ws.Range("E16").Value = "1901"
$ws.Range("F16").Value = 6625
$ws.Range("G16").Value = 828116

# Remove the two rows that are no longer needed (old JULIO CESAR row and
# the old LEIDIS row whose data was moved up into row 16).
$ws.Range("A17:A18").EntireRow.Delete()

# Refresh the summary totals at the top of the sheet.
$ws.Range("E11").Value = 6625
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1
